$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1542.6364
$ws.Range("I28").Value = 1273.6666
$ws.Range("K28").Value = 1273.6666
$ws.Range("M28").Value = -788.6666
$ws.Range("H43").Value = 3400.1667
$ws.Range("I43").Value = 3400.1667
$ws.Range("K43").Value = 3400.1667
$ws.Range("M43").Value = -3331.1667
$ws.Range("H70").Value = 5017.625
$ws.Range("I70").Value = 3419.889
$ws.Range("J70").Value = 7071.857
$ws.Range("K70").Value = 10259.667
$ws.Range("L70").Value = 21215.571
$ws.Range("M70").Value = -9989.667000000001
$ws.Range("N70").Value = -21755.571
$ws.Range("H73").Value = 5017.625
$ws.Range("I73").Value = 3419.889
$ws.Range("J73").Value = 7071.857
$ws.Range("K73").Value = 10259.667
$ws.Range("L73").Value = 21215.571
$ws.Range("M73").Value = -9323.667000000001
$ws.Range("N73").Value = -23087.571
$ws.Range("H80").Value = 5481
$ws.Range("J80").Value = 4524.5
$ws.Range("L80").Value = 13573.5
$ws.Range("N80").Value = -15569.5
$ws.Range("H83").Value = 5481
$ws.Range("J83").Value = 4524.5
$ws.Range("L83").Value = 40720.5
$ws.Range("N83").Value = -50704.5
$ws.Range("H96").Value = 1208.9166
$ws.Range("I96").Value = 1310.4546
$ws.Range("K96").Value = 3931.3638
$ws.Range("M96").Value = -2558.3638
$ws.Range("H100").Value = 2536.3845
$ws.Range("I100").Value = 2370.3635
$ws.Range("K100").Value = 2370.3635
$ws.Range("M100").Value = -1829.3635
$ws.Range("H137").Value = 2030.44
$ws.Range("I137").Value = 1581.7142
$ws.Range("K137").Value = 4745.142599999999
$ws.Range("M137").Value = -2195.142599999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 574
$ws.Range("I5").Value = 680
$ws.Range("J5").Value = 44
$ws.Range("K5").Value = 680
$ws.Range("L5").Value = 44
$ws.Range("M5").Value = -568
$ws.Range("N5").Value = -268
$ws.Range("H32").Value = 3937.484
$ws.Range("I32").Value = 3005.6206
$ws.Range("K32").Value = 3005.6206
$ws.Range("M32").Value = -2718.6206
$ws.Range("H122").Value = 1452.9354
$ws.Range("I122").Value = 1512.138
$ws.Range("K122").Value = 4536.414
$ws.Range("M122").Value = -2086.414
$ws.Range("H132").Value = 9705.714
$ws.Range("I132").Value = 9656.666999999999
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 28970.001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -26440.001
$ws.Range("N132").Value = -35060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 574
$ws.Range("I4").Value = 680
$ws.Range("J4").Value = 44
$ws.Range("K4").Value = 680
$ws.Range("L4").Value = 44
$ws.Range("M4").Value = -565
$ws.Range("N4").Value = -274
$ws.Range("H64").Value = 1992.4166
$ws.Range("I64").Value = 1103
$ws.Range("J64").Value = 2627.7144
$ws.Range("K64").Value = 1103
$ws.Range("L64").Value = 2627.7144
$ws.Range("M64").Value = -878
$ws.Range("N64").Value = -3077.7144
$ws.Range("H67").Value = 1992.4166
$ws.Range("I67").Value = 1103
$ws.Range("J67").Value = 2627.7144
$ws.Range("K67").Value = 1103
$ws.Range("L67").Value = 2627.7144
$ws.Range("M67").Value = -323
$ws.Range("N67").Value = -4187.7144
$ws.Range("H114").Value = 49995
$ws.Range("J114").Value = 49995
$ws.Range("L114").Value = 49995
$ws.Range("N114").Value = -58673

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 375.63635
$ws.Range("I7").Value = 435.875
$ws.Range("J7").Value = 215
$ws.Range("K7").Value = 435.875
$ws.Range("L7").Value = 215
$ws.Range("M7").Value = -322.875
$ws.Range("N7").Value = -441
$ws.Range("H22").Value = 10000493
$ws.Range("I22").Value = 657
$ws.Range("J22").Value = 40000000
$ws.Range("K22").Value = 657
$ws.Range("L22").Value = 40000000
$ws.Range("M22").Value = -307
$ws.Range("N22").Value = -40000700
$ws.Range("H31").Value = 1871.1111
$ws.Range("I31").Value = 1855.25
$ws.Range("K31").Value = 1855.25
$ws.Range("M31").Value = -1560.25
$ws.Range("H34").Value = 1871.1111
$ws.Range("I34").Value = 1855.25
$ws.Range("K34").Value = 1855.25
$ws.Range("M34").Value = -1653.25
$ws.Range("H86").Value = 7934.625
$ws.Range("I86").Value = 7492.3335
$ws.Range("K86").Value = 7492.3335
$ws.Range("M86").Value = -6369.3335
$ws.Range("H89").Value = 7934.625
$ws.Range("I89").Value = 7492.3335
$ws.Range("K89").Value = 37461.6675
$ws.Range("M89").Value = -31845.6675
$ws.Range("H107").Value = 1197.9166
$ws.Range("J107").Value = 2020.8
$ws.Range("L107").Value = 2020.8
$ws.Range("N107").Value = -5860.8
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6470
$ws.Range("H134").Value = 3374.4
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1318403.9
$ws.Range("I4").Value = 17546.47
$ws.Range("K4").Value = 52639.41
$ws.Range("M4").Value = -52527.41
$ws.Range("H62").Value = 10600
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 10600
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H98").Value = 896.6667
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 1095
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 3285
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -6281
$ws.Range("H131").Value = 66837.586
$ws.Range("I131").Value = 1628
$ws.Range("J131").Value = 132047.17
$ws.Range("K131").Value = 4884
$ws.Range("L131").Value = 396141.51
$ws.Range("M131").Value = 156
$ws.Range("N131").Value = -406221.51

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 88.2
$ws.Range("J2").Value = 10
$ws.Range("L2").Value = 10
$ws.Range("N2").Value = -236
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 669
$ws.Range("I9").Value = 669
$ws.Range("K9").Value = 669
$ws.Range("M9").Value = -445
$ws.Range("H22").Value = 943.625
$ws.Range("I22").Value = 591.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 591.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -296.5
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 943.625
$ws.Range("I27").Value = 591.5
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 591.5
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -484.5
$ws.Range("N27").Value = -2214
$ws.Range("H32").Value = 12991.5
$ws.Range("I32").Value = 12991.5
$ws.Range("K32").Value = 12991.5
$ws.Range("M32").Value = -12674.5
$ws.Range("H40").Value = 4904
$ws.Range("I40").Value = 4904
$ws.Range("K40").Value = 4904
$ws.Range("M40").Value = -4768
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("H132").Value = 4838.4
$ws.Range("I132").Value = 4173
$ws.Range("K132").Value = 12519
$ws.Range("M132").Value = -9989

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3257.5217
$ws.Range("I132").Value = 2890.6843
$ws.Range("K132").Value = 8672.052899999999
$ws.Range("M132").Value = -6142.052899999999
